$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Update the cached "datetimeFigureOut" footer date (22.03.2012 -> 23.03.2012)
#    on the slide master and on every slide layout.
# ---------------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "22.03.2012") {
                $tr.Text = "23.03.2012"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    Update-DatePlaceholder $layout.Shapes
}

# ---------------------------------------------------------------------------
# 2. Slide 6, shape "Textfeld 117": widen the box and reword
#    "Mapping depends on " -> "Mapping depending on " (split into new runs).
# ---------------------------------------------------------------------------
$slide6 = $p.Slides.Item(6)
$shp = $slide6.Shapes.Item("Textfeld 117")

# Widen the textbox (1205779 EMU -> 1348446 EMU), leave position/height as-is.
$shp.Width = 106.17686

$tr = $shp.TextFrame.TextRange

# "depends" -> "depending"
$full = $tr.Text
$idx = $full.IndexOf("depends")
$run = $tr.Characters($idx + 1, 7)
$run.Text = "depending"

# Split the following " on " run into " " + "on " (same visible text,
# matches the run layout produced by the original edit).
$full2 = $tr.Text
$idx2 = $full2.IndexOf(" on ")
$spaceRun = $tr.Characters($idx2 + 1, 1)
$spaceRun.Text = " "
$onRun = $tr.Characters($idx2 + 2, 3)
$onRun.Text = "on "
